# Append: 2025-09-27 12:32 JST
# Refresh the "ランサーズ" sheet with the latest scrape: existing rows 2-3 keep
# their content (only the fetch timestamp moves forward), a few rows are
# re-shuffled/replaced, and three brand-new listings are appended, growing the
# sheet from A1:H7 to A1:H10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-09-27 12:32:16"

# r, title, category, price, deadline, url, score, skill-summary (blank = none)
$rows = @(
  @(2,  "【限定タスク】SIM AI の Google 認証ログイン機能の「最終調整」のみ代行(環境構築済み)", "システム開発", "20,000 円 ~ 50,000 円 / 固定",  "期限情報なし", "https://www.lancers.jp/work/detail/5401688", 298, "🔥AI,Ai"),
  @(3,  "SIM AI の Google 認証ログイン機能の「最終調整」のみ代行(環境構築済み)",               "システム開発", "10,000 円 ~ 20,000 円 / 固定",  "期限情報なし", "https://www.lancers.jp/work/detail/5401604", 295, "🔥AI,Ai"),
  @(4,  "NOTEへのいいね・フォローを行うツールの販売",                                            "システム開発", "20,000 円 ~ 50,000 円 / 固定",  "期限情報なし", "https://www.lancers.jp/work/detail/5401806", 68,  "◆ツール"),
  @(5,  "iOSアプリのAdMobメディエーション入札接続とeCPM改善",                                    "システム開発", "100,000 円 ~ 200,000 円 / 固定","期限情報なし", "https://www.lancers.jp/work/detail/5401602", 38,  "◇アプリ"),
  @(6,  "初回 シティヘブンの自動予約システム",                                                    "システム開発", "20,000 円 ~ 50,000 円 / 固定",  "期限情報なし", "https://www.lancers.jp/work/detail/5401800", 28,  ""),
  @(7,  "限定公開 PR 限定公開の仕事",                                                             "システム開発", "20,000 円 ~ 50,000 円 / 固定",  "期限情報なし", "https://www.lancers.jp/work/detail/5399347", 13,  ""),
  @(8,  "【急募】AWS環境構築できる方!",                                                           "システム開発", "20,000 円 ~ 50,000 円 / 固定",  "期限情報なし", "https://www.lancers.jp/work/detail/5401736", 13,  ""),
  @(9,  "Android kotlinの画像ファイル 拡張子取得",                                                "システム開発", "~ 5,000 円 / 固定",             "期限情報なし", "https://www.lancers.jp/work/detail/5401572", 10,  ""),
  @(10, "【急募】ワードプレスのローカル環境構築をサポートしてくれる方",                            "システム開発", "~ 5,000 円 / 固定",             "期限情報なし", "https://www.lancers.jp/work/detail/5401534", 10,  "")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]

    if ($row[7] -ne "") {
        $ws.Cells.Item($r, 8).Value = $row[7]
    } else {
        $ws.Cells.Item($r, 8).ClearContents()
    }
}

# The hyperlink collection doesn't auto-shift with cell writes, and per-item
# Delete() is unreliable here, so rebuild it wholesale: wipe every hyperlink
# on the sheet, then re-add one per data row against the just-written URLs.
$ws.Hyperlinks.Delete()
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[5])
}
